{"js": "// Apply the documented text corrections throughout the report body.\n// Each fix is performed via a unique, narrow search string so we only\n// ever touch the intended occurrence.\n\nasync function replaceOnce(body, searchText, replacement, options) {\n  const results = body.search(searchText, options || { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\n      \"Expected exactly one match for '\" +\n        searchText +\n        \"' but found \" +\n        results.items.length\n    );\n  }\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst body = context.document.body;\n\n// 1) \" With in this \" -> \" Within this \"\nawait replaceOnce(body, \"With in this\", \"Within this\");\n\n// 2) \"pygmame\" -> \"pygame\"\nawait replaceOnce(body, \"pygmame\", \"pygame\");\n\n// 3) \"ailien_invasion\" -> \"alien_invasion\"\nawait replaceOnce(body, \"ailien_invasion\", \"alien_invasion\");\n\n// 4) \"has its perimeter defined we utilize\" -> \"... defined, we utilize\"\nawait replaceOnce(\n  body,\n  \"has its perimeter defined we utilize\",\n  \"has its perimeter defined, we utilize\"\n);\n\n// 5) \"although very short allows us\" -> \"... short, allows us\"\nawait replaceOnce(\n  body,\n  \"although very short allows us\",\n  \"although very short, allows us\"\n);\n\n// 6) \"players input\" -> \"player\\u2019s input\"\nawait replaceOnce(body, \"players input\", \"player\\u2019s input\");\n\n// 7) \"wavetops of project, Tim created\" -> \"... project, Adam created\"\nawait replaceOnce(\n  body,\n  \"wavetops of project, Tim created\",\n  \"wavetops of project, Adam created\"\n);\n\n// 8) \"and adjust scoring\" -> \"and adjusting scoring\"\nawait replaceOnce(body, \"and adjust scoring\", \"and adjusting scoring\");\n", "ps1": "# Apply the documented text corrections throughout the report body.\n# Each fix uses Find/Replace with a unique, narrow search string so we\n# only ever touch the intended occurrence.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Once {\n    param(\n        [string]$FindText,\n        [string]$ReplaceText,\n        [bool]$MatchCase = $true\n    )\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($FindText, $MatchCase, $false, $false, $false, $false, $true, 1, $false, $ReplaceText, 2)\n}\n\n# 1) \" With in this \" -> \" Within this \"\nReplace-Once \"With in this\" \"Within this\"\n\n# 2) \"pygmame\" -> \"pygame\"\nReplace-Once \"pygmame\" \"pygame\"\n\n# 3) \"ailien_invasion\" -> \"alien_invasion\"\nReplace-Once \"ailien_invasion\" \"alien_invasion\"\n\n# 4) \"has its perimeter defined we utilize\" -> \"... defined, we utilize\"\nReplace-Once \"has its perimeter defined we utilize\" \"has its perimeter defined, we utilize\"\n\n# 5) \"although very short allows us\" -> \"... short, allows us\"\nReplace-Once \"although very short allows us\" \"although very short, allows us\"\n\n# 6) \"players input\" -> \"player's input\" (curly apostrophe)\nReplace-Once \"players input\" \"player\u2019s input\"\n\n# 7) \"wavetops of project, Tim created\" -> \"... project, Adam created\"\nReplace-Once \"wavetops of project, Tim created\" \"wavetops of project, Adam created\"\n\n# 8) \"and adjust scoring\" -> \"and adjusting scoring\"\nReplace-Once \"and adjust scoring\" \"and adjusting scoring\"\n\nWrite-Output \"done\"\n"}
